{"js": "// Map of original text -> new text, in document order.\n// (1 date heading + 25 division problems)\nconst replacements = [\n  [\"2026-01-04 Sunday\", \"2026-01-05 Monday\"],\n  [\"24\u00f77=\", \"90\u00f74=\"],\n  [\"48\u00f73=\", \"79\u00f72=\"],\n  [\"97\u00f79=\", \"19\u00f77=\"],\n  [\"96\u00f76=\", \"75\u00f79=\"],\n  [\"34\u00f76=\", \"98\u00f73=\"],\n  [\"28\u00f79=\", \"39\u00f72=\"],\n  [\"21\u00f76=\", \"41\u00f76=\"],\n  [\"40\u00f74=\", \"41\u00f72=\"],\n  [\"77\u00f76=\", \"68\u00f75=\"],\n  [\"92\u00f74=\", \"79\u00f77=\"],\n  [\"48\u00f76=\", \"82\u00f78=\"],\n  [\"67\u00f73=\", \"23\u00f79=\"],\n  [\"80\u00f76=\", \"74\u00f76=\"],\n  [\"71\u00f73=\", \"96\u00f72=\"],\n  [\"71\u00f75=\", \"34\u00f73=\"],\n  [\"29\u00f72=\", \"13\u00f74=\"],\n  [\"18\u00f72=\", \"56\u00f79=\"],\n  [\"57\u00f76=\", \"13\u00f76=\"],\n  [\"90\u00f75=\", \"19\u00f73=\"],\n  [\"64\u00f77=\", \"87\u00f77=\"],\n  [\"44\u00f75=\", \"80\u00f78=\"],\n  [\"38\u00f72=\", \"91\u00f74=\"],\n  [\"47\u00f77=\", \"43\u00f79=\"],\n  [\"26\u00f75=\", \"91\u00f72=\"],\n  [\"59\u00f72=\", \"16\u00f78=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Map of original text -> new text, in document order.\n# (1 date heading + 25 division problems)\n$replacements = @(\n    @(\"2026-01-04 Sunday\", \"2026-01-05 Monday\"),\n    @(\"24\u00f77=\", \"90\u00f74=\"),\n    @(\"48\u00f73=\", \"79\u00f72=\"),\n    @(\"97\u00f79=\", \"19\u00f77=\"),\n    @(\"96\u00f76=\", \"75\u00f79=\"),\n    @(\"34\u00f76=\", \"98\u00f73=\"),\n    @(\"28\u00f79=\", \"39\u00f72=\"),\n    @(\"21\u00f76=\", \"41\u00f76=\"),\n    @(\"40\u00f74=\", \"41\u00f72=\"),\n    @(\"77\u00f76=\", \"68\u00f75=\"),\n    @(\"92\u00f74=\", \"79\u00f77=\"),\n    @(\"48\u00f76=\", \"82\u00f78=\"),\n    @(\"67\u00f73=\", \"23\u00f79=\"),\n    @(\"80\u00f76=\", \"74\u00f76=\"),\n    @(\"71\u00f73=\", \"96\u00f72=\"),\n    @(\"71\u00f75=\", \"34\u00f73=\"),\n    @(\"29\u00f72=\", \"13\u00f74=\"),\n    @(\"18\u00f72=\", \"56\u00f79=\"),\n    @(\"57\u00f76=\", \"13\u00f76=\"),\n    @(\"90\u00f75=\", \"19\u00f73=\"),\n    @(\"64\u00f77=\", \"87\u00f77=\"),\n    @(\"44\u00f75=\", \"80\u00f78=\"),\n    @(\"38\u00f72=\", \"91\u00f74=\"),\n    @(\"47\u00f77=\", \"43\u00f79=\"),\n    @(\"26\u00f75=\", \"91\u00f72=\"),\n    @(\"59\u00f72=\", \"16\u00f78=\")\n)\n\n$d = $word.ActiveDocument\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Forward = $true\n    $find.Wrap = 1  # wdFindContinue\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n\n    $find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2) | Out-Null\n}\n"}
